$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4001.026
$ws.Range("I64").Value = 3943.3076
$ws.Range("J64").Value = 4060.2632
$ws.Range("K64").Value = 3943.3076
$ws.Range("L64").Value = 4060.2632
$ws.Range("M64").Value = -3695.3076
$ws.Range("N64").Value = -4556.263199999999
$ws.Range("H67").Value = 4001.026
$ws.Range("I67").Value = 3943.3076
$ws.Range("J67").Value = 4060.2632
$ws.Range("K67").Value = 3943.3076
$ws.Range("L67").Value = 4060.2632
$ws.Range("M67").Value = -3085.3076
$ws.Range("N67").Value = -5776.263199999999
$ws.Range("H92").Value = 4114.6665
$ws.Range("I92").Value = 172
$ws.Range("J92").Value = 12000
$ws.Range("K92").Value = 172
$ws.Range("L92").Value = 12000
$ws.Range("M92").Value = 1076
$ws.Range("N92").Value = -14496
$ws.Range("H99").Value = 622
$ws.Range("I99").Value = 237.75
$ws.Range("J99").Value = 1390.5
$ws.Range("K99").Value = 713.25
$ws.Range("L99").Value = 4171.5
$ws.Range("M99").Value = 784.75
$ws.Range("N99").Value = -7167.5
$ws.Range("H113").Value = 1456.3914
$ws.Range("I113").Value = 949.875
$ws.Range("J113").Value = 1726.5333
$ws.Range("K113").Value = 949.875
$ws.Range("L113").Value = 1726.5333
$ws.Range("M113").Value = 2304.125
$ws.Range("N113").Value = -8234.533299999999
$ws.Range("H125").Value = 895
$ws.Range("I125").Value = 895
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8055
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -5595
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 719.7727
$ws.Range("I127").Value = 499.1875
$ws.Range("J127").Value = 1308
$ws.Range("K127").Value = 1497.5625
$ws.Range("L127").Value = 3924
$ws.Range("M127").Value = 3462.4375
$ws.Range("N127").Value = -13844
$ws.Range("H131").Value = 2561
$ws.Range("I131").Value = 1051.5834
$ws.Range("J131").Value = 5579.8335
$ws.Range("K131").Value = 3154.7502
$ws.Range("L131").Value = 16739.5005
$ws.Range("M131").Value = 1885.2498
$ws.Range("N131").Value = -26819.5005
$ws.Range("H135").Value = 1502
$ws.Range("I135").Value = 1144
$ws.Range("J135").Value = 1860
$ws.Range("K135").Value = 10296
$ws.Range("L135").Value = 16740
$ws.Range("M135").Value = -7761
$ws.Range("N135").Value = -21810
$ws.Range("H137").Value = 27753.105
$ws.Range("I137").Value = 1310.8966
$ws.Range("J137").Value = 112955.78
$ws.Range("K137").Value = 3932.6898
$ws.Range("L137").Value = 338867.34
$ws.Range("M137").Value = -1382.6898
$ws.Range("N137").Value = -343967.34
$ws.Range("H138").Value = 2220.6558
$ws.Range("I138").Value = 1279.8975
$ws.Range("J138").Value = 3888.3635
$ws.Range("K138").Value = 3839.6925
$ws.Range("L138").Value = 11665.0905
$ws.Range("M138").Value = 1300.3075
$ws.Range("N138").Value = -21945.0905
$ws.Range("H141").Value = 1508.5883
$ws.Range("I141").Value = 928.8333
$ws.Range("J141").Value = 2900
$ws.Range("K141").Value = 2786.4999
$ws.Range("L141").Value = 8700
$ws.Range("M141").Value = 2393.5001
$ws.Range("N141").Value = -19060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3867.4
$ws.Range("I61").Value = 3808.5715
$ws.Range("J61").Value = 4004.6667
$ws.Range("K61").Value = 3808.5715
$ws.Range("L61").Value = 4004.6667
$ws.Range("M61").Value = -3596.5715
$ws.Range("N61").Value = -4428.6667
$ws.Range("H74").Value = 5275.375
$ws.Range("I74").Value = 6305
$ws.Range("J74").Value = 1362.8
$ws.Range("K74").Value = 6305
$ws.Range("L74").Value = 1362.8
$ws.Range("M74").Value = -5431
$ws.Range("N74").Value = -3110.8
$ws.Range("H77").Value = 5275.375
$ws.Range("I77").Value = 6305
$ws.Range("J77").Value = 1362.8
$ws.Range("K77").Value = 31525
$ws.Range("L77").Value = 6814
$ws.Range("M77").Value = -27157
$ws.Range("N77").Value = -15550
$ws.Range("H97").Value = 522.13043
$ws.Range("I97").Value = 493.93332
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 493.93332
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = 2.066680000000019
$ws.Range("N97").Value = -1567
$ws.Range("H102").Value = 1767.1333
$ws.Range("I102").Value = 1309.375
$ws.Range("J102").Value = 2290.2856
$ws.Range("K102").Value = 1309.375
$ws.Range("L102").Value = 2290.2856
$ws.Range("M102").Value = 312.625
$ws.Range("N102").Value = -5534.2856
$ws.Range("H122").Value = 2264192.2
$ws.Range("I122").Value = 2942640
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 8827920
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -8825470
$ws.Range("N122").Value = -13000
$ws.Range("H132").Value = 4430.9
$ws.Range("I132").Value = 4299.5
$ws.Range("J132").Value = 4956.5
$ws.Range("K132").Value = 12898.5
$ws.Range("L132").Value = 14869.5
$ws.Range("M132").Value = -10368.5
$ws.Range("N132").Value = -19929.5
$ws.Range("H136").Value = 3867.4
$ws.Range("I136").Value = 3808.5715
$ws.Range("J136").Value = 4004.6667
$ws.Range("K136").Value = 11425.7145
$ws.Range("L136").Value = 12014.0001
$ws.Range("M136").Value = -8875.7145
$ws.Range("N136").Value = -17114.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1913
$ws.Range("I107").Value = 1932.7693
$ws.Range("J107").Value = 1827.3334
$ws.Range("K107").Value = 1932.7693
$ws.Range("L107").Value = 1827.3334
$ws.Range("M107").Value = -12.76929999999993
$ws.Range("N107").Value = -5667.3334
$ws.Range("H134").Value = 1223.1608
$ws.Range("I134").Value = 1253
$ws.Range("J134").Value = 918.8
$ws.Range("K134").Value = 3759
$ws.Range("L134").Value = 2756.4
$ws.Range("M134").Value = -1224
$ws.Range("N134").Value = -7826.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3898.9697
$ws.Range("I31").Value = 3566.7856
$ws.Range("J31").Value = 5759.2
$ws.Range("K31").Value = 3566.7856
$ws.Range("L31").Value = 5759.2
$ws.Range("M31").Value = -3271.7856
$ws.Range("N31").Value = -6349.2
$ws.Range("H34").Value = 3898.9697
$ws.Range("I34").Value = 3566.7856
$ws.Range("J34").Value = 5759.2
$ws.Range("K34").Value = 3566.7856
$ws.Range("L34").Value = 5759.2
$ws.Range("M34").Value = -3364.7856
$ws.Range("N34").Value = -6163.2
$ws.Range("H107").Value = 500.56412
$ws.Range("I107").Value = 426.2
$ws.Range("J107").Value = 633.3570999999999
$ws.Range("K107").Value = 426.2
$ws.Range("L107").Value = 633.3570999999999
$ws.Range("M107").Value = 1493.8
$ws.Range("N107").Value = -4473.3571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 807.3570999999999
$ws.Range("I122").Value = 664
$ws.Range("J122").Value = 1333
$ws.Range("K122").Value = 5976
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -3526
$ws.Range("N122").Value = -16897
$ws.Range("H131").Value = 748.5
$ws.Range("I131").Value = 356
$ws.Range("J131").Value = 1028.8572
$ws.Range("K131").Value = 1068
$ws.Range("L131").Value = 3086.5716
$ws.Range("M131").Value = 3972
$ws.Range("N131").Value = -13166.5716
$ws.Range("H132").Value = 722750.8
$ws.Range("I132").Value = 743
$ws.Range("J132").Value = 1444758.6
$ws.Range("K132").Value = 6687
$ws.Range("L132").Value = 13002827.4
$ws.Range("M132").Value = -4157
$ws.Range("N132").Value = -13007887.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 269.72223
$ws.Range("I107").Value = 229.07692
$ws.Range("J107").Value = 375.4
$ws.Range("K107").Value = 229.07692
$ws.Range("L107").Value = 375.4
$ws.Range("M107").Value = 1690.92308
$ws.Range("N107").Value = -4215.4
$ws.Range("H123").Value = 18609.736
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 18609.736
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 18609.736
$ws.Range("N123").Value = -23509.736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1790.3636
$ws.Range("I122").Value = 1823
$ws.Range("J122").Value = 1703.3334
$ws.Range("K122").Value = 5469
$ws.Range("L122").Value = 5110.0002
$ws.Range("M122").Value = -3019
$ws.Range("N122").Value = -10010.0002
$ws.Range("H132").Value = 4190.227
$ws.Range("I132").Value = 5199.364
$ws.Range("J132").Value = 3181.0908
$ws.Range("K132").Value = 15598.092
$ws.Range("L132").Value = 9543.2724
$ws.Range("M132").Value = -13068.092
$ws.Range("N132").Value = -14603.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 11153.846
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 15000
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = -9712
$ws.Range("N15").Value = -15576
$ws.Range("H122").Value = 2064.1667
$ws.Range("I122").Value = 1535.8422
$ws.Range("J122").Value = 2654.647
$ws.Range("K122").Value = 4607.5266
$ws.Range("L122").Value = 7963.941
$ws.Range("M122").Value = -2157.5266
$ws.Range("N122").Value = -12863.941
$ws.Range("H132").Value = 1809.85
$ws.Range("I132").Value = 1129.2941
$ws.Range("J132").Value = 5666.3335
$ws.Range("K132").Value = 3387.8823
$ws.Range("L132").Value = 16999.0005
$ws.Range("M132").Value = -857.8823000000002
$ws.Range("N132").Value = -22059.0005
